$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cosmetic: remembered on-screen window position (best effort; not core data).
$win = $excel.ActiveWindow
$win.Left = 1275
$win.Top = 2895

# A1 was a numeric 1; change it to a quote-prefixed text value "1" (stored as a shared
# string but formatted as General, the leading apostrophe forces text entry) and
# center-align it, same as the newly added A2:A10 cells below.
$ws.Range("A1:A10").HorizontalAlignment = -4108  # xlCenter

for ($i = 1; $i -le 10; $i++) {
    $cell = $ws.Cells.Item($i, 1)
    $cell.Value = "'" + [string]$i
}

[void]$ws.Range("B14").Select()
